# Auto-generated edit script applying numeric updates described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50924

$ws.Range("H64").Value = 2942.8572
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752

$ws.Range("H67").Value = 2942.8572
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142

$ws.Range("H112").Value = 1431.49
$ws.Range("J112").Value = 1431.8182
$ws.Range("L112").Value = 4295.4546
$ws.Range("N112").Value = -6511.4546

$ws.Range("H116").Value = 723236.5600000001
$ws.Range("I116").Value = 1113922.2
$ws.Range("K116").Value = 1113922.2
$ws.Range("M116").Value = -1110480.2

$ws.Range("H127").Value = 1737.6
$ws.Range("I127").Value = 747.3333
$ws.Range("K127").Value = 2241.9999
$ws.Range("M127").Value = 2718.0001

$ws.Range("H129").Value = 912.8889
$ws.Range("I129").Value = 399
$ws.Range("K129").Value = 1197
$ws.Range("M129").Value = 3803

$ws.Range("H132").Value = 99701.414
$ws.Range("I132").Value = 122820.55
$ws.Range("K132").Value = 368461.65
$ws.Range("M132").Value = -365931.65

$ws.Range("H137").Value = 2248.4888
$ws.Range("I137").Value = 1521.6757
$ws.Range("J137").Value = 5610
$ws.Range("K137").Value = 4565.0271
$ws.Range("L137").Value = 16830
$ws.Range("M137").Value = -2015.0271
$ws.Range("N137").Value = -21930

$ws.Range("H138").Value = 2701.0942
$ws.Range("I138").Value = 1628.7
$ws.Range("J138").Value = 3031.0615
$ws.Range("K138").Value = 4886.1
$ws.Range("L138").Value = 9093.184499999999
$ws.Range("M138").Value = 253.8999999999996
$ws.Range("N138").Value = -19373.1845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15366.706
$ws.Range("I32").Value = 12892.818
$ws.Range("K32").Value = 12892.818
$ws.Range("M32").Value = -12605.818

$ws.Range("H45").Value = 2032.8823
$ws.Range("I45").Value = 1100
$ws.Range("J45").Value = 2541.7273
$ws.Range("K45").Value = 1100
$ws.Range("L45").Value = 2541.7273
$ws.Range("M45").Value = -723
$ws.Range("N45").Value = -3295.7273

$ws.Range("H61").Value = 1269.4103
$ws.Range("I61").Value = 871.76
$ws.Range("J61").Value = 1979.5
$ws.Range("K61").Value = 871.76
$ws.Range("L61").Value = 1979.5
$ws.Range("M61").Value = -659.76
$ws.Range("N61").Value = -2403.5

$ws.Range("H63").Value = 13854341
$ws.Range("I63").Value = 19789688
$ws.Range("K63").Value = 19789688
$ws.Range("M63").Value = -19789002

$ws.Range("H66").Value = 13854341
$ws.Range("I66").Value = 19789688
$ws.Range("K66").Value = 98948440
$ws.Range("M66").Value = -98945008

$ws.Range("H132").Value = 2717.697
$ws.Range("I132").Value = 1275.3
$ws.Range("J132").Value = 4936.769
$ws.Range("K132").Value = 3825.9
$ws.Range("L132").Value = 14810.307
$ws.Range("M132").Value = -1295.9
$ws.Range("N132").Value = -19870.307

$ws.Range("H136").Value = 1269.4103
$ws.Range("I136").Value = 871.76
$ws.Range("J136").Value = 1979.5
$ws.Range("K136").Value = 2615.28
$ws.Range("L136").Value = 5938.5
$ws.Range("M136").Value = -65.27999999999975
$ws.Range("N136").Value = -11038.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2764.3
$ws.Range("I99").Value = 1114.6471
$ws.Range("J99").Value = 4921.5386
$ws.Range("K99").Value = 1114.6471
$ws.Range("L99").Value = 4921.5386
$ws.Range("M99").Value = 383.3529000000001
$ws.Range("N99").Value = -7917.5386

$ws.Range("H107").Value = 1504
$ws.Range("I107").Value = 1402.2
$ws.Range("J107").Value = 2013
$ws.Range("K107").Value = 1402.2
$ws.Range("L107").Value = 2013
$ws.Range("M107").Value = 517.8
$ws.Range("N107").Value = -5853

$ws.Range("H134").Value = 2601.869
$ws.Range("I134").Value = 1572.6511
$ws.Range("J134").Value = 5060.5557
$ws.Range("K134").Value = 4717.9533
$ws.Range("L134").Value = 15181.6671
$ws.Range("M134").Value = -2182.9533
$ws.Range("N134").Value = -20251.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2077.2
$ws.Range("I105").Value = 1531.1111
$ws.Range("J105").Value = 2896.3333
$ws.Range("K105").Value = 1531.1111
$ws.Range("L105").Value = 2896.3333
$ws.Range("M105").Value = 215.8888999999999
$ws.Range("N105").Value = -6390.3333

$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676

$ws.Range("H122").Value = 1942.8
$ws.Range("I122").Value = 1441.4584
$ws.Range("J122").Value = 3948.1667
$ws.Range("K122").Value = 4324.3752
$ws.Range("L122").Value = 11844.5001
$ws.Range("M122").Value = -1874.3752
$ws.Range("N122").Value = -16744.5001

$ws.Range("H132").Value = 2818.08
$ws.Range("I132").Value = 1726.8462
$ws.Range("J132").Value = 4000.25
$ws.Range("K132").Value = 5180.5386
$ws.Range("L132").Value = 12000.75
$ws.Range("M132").Value = -2650.5386
$ws.Range("N132").Value = -17060.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 510.2143
$ws.Range("I40").Value = 142.875
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 571.5
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -502.5
$ws.Range("N40").Value = -4138

$ws.Range("H62").Value = 8188
$ws.Range("J62").Value = 9985
$ws.Range("L62").Value = 29955
$ws.Range("N62").Value = -31327

$ws.Range("H65").Value = 8188
$ws.Range("J65").Value = 9985
$ws.Range("L65").Value = 89865
$ws.Range("N65").Value = -96729

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990

$ws.Range("H102").Value = 3064.0952
$ws.Range("I102").Value = 2419.3333
$ws.Range("J102").Value = 6932.6665
$ws.Range("K102").Value = 2419.3333
$ws.Range("L102").Value = 6932.6665
$ws.Range("M102").Value = -797.3332999999998
$ws.Range("N102").Value = -10176.6665

$ws.Range("H120").Value = 39766.668
$ws.Range("J120").Value = 39766.668
$ws.Range("L120").Value = 39766.668
$ws.Range("N120").Value = -49442.668

$ws.Range("H126").Value = 4103.2183
$ws.Range("I126").Value = 2951.0637
$ws.Range("J126").Value = 5457
$ws.Range("K126").Value = 8853.1911
$ws.Range("L126").Value = 16371
$ws.Range("M126").Value = -6383.1911
$ws.Range("N126").Value = -21311

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3277
$ws.Range("J136").Value = 4917.2666
$ws.Range("L136").Value = 14751.7998
$ws.Range("N136").Value = -19851.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H121").Value = 27944.834
$ws.Range("J121").Value = 27944.834
$ws.Range("L121").Value = 27944.834
$ws.Range("N121").Value = -31438.834

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
